$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells receiving a new numeric-looking price value need to be forced to
# Text format first, otherwise Excel silently re-types the string as a
# number (dropping trailing zeros / switching to scientific notation).

# Row 2
$ws.Range("D2").Value = "74.725.84"
$ws.Range("E2").Value = "  +0.29%  "

# Row 3
$ws.Range("D3").Value = "2.838.56"
$ws.Range("E3").Value = "  +9.53%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.22"
$ws.Range("E5").Value = "  +3.88%  "

# Row 6
$ws.Range("B6").Value = "Solana"
$ws.Range("C6").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "188.83"
$ws.Range("E6").Value = "  +1.51%  "

# Row 7
$ws.Range("E7").Value = "  -0.07%  "

# Row 8
$ws.Range("E8").Value = "  +3.92%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.194"
$ws.Range("E9").Value = "  -6.22%  "

# Row 10
$ws.Range("D10").Value = "2.838.83"
$ws.Range("E10").Value = "  +9.63%  "

# Row 11
$ws.Range("E11").Value = "  -0.24%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.372"
$ws.Range("E12").Value = "  +3.56%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.89"
$ws.Range("E13").Value = "  +2.03%  "

# Row 14
$ws.Range("D14").Value = "3.363.92"
$ws.Range("E14").Value = "  +9.67%  "

# Row 15
$ws.Range("D15").Value = "74.788.37"
$ws.Range("E15").Value = "  +0.54%  "

# Row 16
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000189"
$ws.Range("E16").Value = "  -1.91%  "

# Row 17
$ws.Range("B17").Value = "Avalanche"
$ws.Range("C17").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.30"
$ws.Range("E17").Value = "  +4.15%  "

# Row 18
$ws.Range("D18").Value = "2.840.72"
$ws.Range("E18").Value = "  +9.91%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.19"
$ws.Range("E19").Value = "  +8.35%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.50"
$ws.Range("E20").Value = "  +6.91%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "376.91"
$ws.Range("E21").Value = "  -0.29%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.28"
$ws.Range("E22").Value = "  -1.19%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.14"
$ws.Range("E23").Value = "  +1.65%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.22"
$ws.Range("E24").Value = "  +0.16%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  -0.03%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "70.79"
$ws.Range("E26").Value = "  +1.35%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.24"
$ws.Range("E27").Value = "  +1.14%  "

# Row 28
$ws.Range("D28").Value = "2.983.50"
$ws.Range("E28").Value = "  +9.82%  "

# Row 29
$ws.Range("E29").Value = "  +4.66%  "

# Row 30
$ws.Range("E30").Value = "  +10.82%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "526.12"
$ws.Range("E32").Value = "  +5.44%  "

# Row 33
$ws.Range("E33").Value = "  +4.53%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.95"
$ws.Range("E34").Value = "  +0.04%  "

# Row 35
$ws.Range("E35").Value = "  +6.13%  "

# Row 36
$ws.Range("E36").Value = "  -0.03%  "

# Row 37
$ws.Range("E37").Value = "  +0.96%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "20.15"
$ws.Range("E38").Value = "  +4.93%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "162.36"
$ws.Range("E39").Value = "  +1.64%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.29"
$ws.Range("E40").Value = "  -0.49%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "185.03"
$ws.Range("E41").Value = "  +24.98%  "

# Row 42
$ws.Range("E42").Value = "  +0.01%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.11"
$ws.Range("E43").Value = "  +2.79%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.341"
$ws.Range("E44").Value = "  +6.83%  "

# Row 45
$ws.Range("E45").Value = "  +1.54%  "

# Row 46
$ws.Range("E46").Value = "  +7.85%  "

# Row 47
$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.38"
$ws.Range("E47").Value = "  -1.84%  "

# Row 48
$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "39.58"
$ws.Range("E48").Value = "  +1.36%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0851"
$ws.Range("E49").Value = "  +4.43%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.574"
$ws.Range("E50").Value = "  +10.52%  "

# Row 51
$ws.Range("E51").Value = "  +3.88%  "
